# Performance Approach.docx edit script
# Summary of the change (per commit message "doc: performance approach done"):
#  - Remove the empty "Computational Efficiency" subsection (it had no body text, just a
#    placeholder "X"), relocating the _GoBack bookmark that used to sit at the end of the
#    "Maintainability" heading to the blank paragraph that now directly follows the intro.
#  - Fill in the placeholder "X" paragraphs under Reliability / Security / Portability /
#    Maintainability with the real evaluation text.
#  - Add two new "improvement" paragraphs after Portability (before Maintainability) and
#    two new paragraphs after Maintainability's body text (before Scalability).
#  - Tidy the Scalability paragraph: remove the spell-check markup around "Github" so the
#    run is one contiguous piece of text.

$d = $word.ActiveDocument

function Get-HeadingParagraph([string]$headingText) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Style.NameLocal -eq "Heading 2" -and $p.Range.Text.TrimEnd([char]13) -eq $headingText) {
            return $p
        }
    }
    return $null
}

# ---------------------------------------------------------------------------
# 1. Relocate the _GoBack bookmark from the Maintainability heading paragraph
#    to the blank paragraph right after the introduction (the paragraph that
#    currently precedes "Computational Efficiency").
# ---------------------------------------------------------------------------
$introBlankPara = $d.Paragraphs.Item(3)
$bookmarkTarget = $introBlankPara.Range.Duplicate
$bookmarkTarget.Collapse(1)  # wdCollapseStart

$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()
$d.Bookmarks.Add("_GoBack", $bookmarkTarget) | Out-Null

# ---------------------------------------------------------------------------
# 2. Remove the whole "Computational Efficiency" subsection: heading + "X" +
#    trailing blank paragraph.
# ---------------------------------------------------------------------------
$compEffHeading = Get-HeadingParagraph("Computational Efficiency")
$compEffBlank = $compEffHeading.Next().Next()
$removeRange = $d.Range($compEffHeading.Range.Start, $compEffBlank.Range.End)
$removeRange.Delete()

# ---------------------------------------------------------------------------
# 3. Fill in the real body text for Reliability / Security / Portability /
#    Maintainability (replacing the placeholder "X").
# ---------------------------------------------------------------------------
$reliabilityHeading = Get-HeadingParagraph("Reliability")
$reliabilityHeading.Next().Range.Text = "The fact that the system will always complete its required job regardless proves that it is reliable, however, as stated later on in this evaluation, it does require the specified files for the program to work properly as it doesn" + [char]0x2019 + "t create its own XML files. However, as the code is modular in its programming style, it can edited and still work in full, for when a part of the program is under construction. "

$securityHeading = Get-HeadingParagraph("Security")
$securityHeading.Next().Range.Text = "The security of the program is easily compromised. It holds user data as static strings in XML files that can be easily opened and viewed using file viewers, including passwords and email combinations. If this was to become a release it would have to be hashed and then decrypted when passwords are entered. This is also an issue due to the fact sensitive plane data is also set up as static information so can also be opened, which could compromise the planes and also be a threat to the passengers safety."

$portabilityHeading = Get-HeadingParagraph("Portability")
$portabilityHeading.Next().Range.Text = "The way the program is set up is it users various XML files that are scattered in the source directory, meaning that if this was moved around different machines, it would have to have all the files added to the same directory otherwise it wouldn" + [char]0x2019 + "t input correctly. It wouldn" + [char]0x2019 + "t let the user know, either, when a file cannot be found and would just stop working, which would take away from the user experience. "

$maintainabilityHeading = Get-HeadingParagraph("Maintainability")
$maintainabilityHeading.Next().Range.Text = "The program can be looked after constantly using the very easily viewable, split up functions and then code can be added and changed. Due to the way the menu is set up, when a menu item is being worked on it can simply be taken out of the menu and not change how the program functions. Obviously, it wouldn" + [char]0x2019 + "t be able to be inputted either. This is a slow way of maintaining a program but it is effective, although basic."

# ---------------------------------------------------------------------------
# 4. Insert the two new paragraphs after the Portability text/blank line,
#    before the Maintainability heading.
# ---------------------------------------------------------------------------
$portabilityBlank = $portabilityHeading.Next().Next()
$insertPoint = $portabilityBlank.Range.Duplicate
$insertPoint.Collapse(0)  # wdCollapseEnd
$insertPoint.InsertParagraphAfter()
$improveTextPara = $insertPoint.Paragraphs.Item(1).Next()
$improveTextPara = $d.Paragraphs.Item($portabilityBlank.Index + 1)
$improveTextPara.Range.Text = "To improve this, we would have to add a directory that holds all the user inf"
$afterFirstRun = $improveTextPara.Range.End - 1
$r2 = $d.Range($afterFirstRun, $afterFirstRun)
$r2.InsertAfter("ormation and the XML files, then this can moved around with the release rather than having to add everything later on. Or, another way of getting around this issue, is to be able to create the XML file if it doesn" + [char]0x2019 + "t exist previously.")

$improveBlankPara = $improveTextPara.Range.Duplicate
$improveBlankPara.Collapse(0)
$improveBlankPara.Start = $improveTextPara.Range.End - 1
$improveBlankPara.End = $improveTextPara.Range.End - 1
$improveBlankPara.InsertParagraphAfter()

# ---------------------------------------------------------------------------
# 5. Insert the two new paragraphs after the Maintainability text/blank line,
#    before the Scalability heading.
# ---------------------------------------------------------------------------
$maintainabilityHeading = Get-HeadingParagraph("Maintainability")
$maintainabilityBlank = $maintainabilityHeading.Next().Next()
$insertPoint2 = $maintainabilityBlank.Range.Duplicate
$insertPoint2.Collapse(0)
$insertPoint2.InsertParagraphAfter()
$issuesPara = $d.Paragraphs.Item($maintainabilityBlank.Index + 1)
$issuesPara.Range.Text = "One of the only big issues with maintaining the program is the large size of it, as it would need someone who knows what they" + [char]0x2019 + "re doing to change the program, as someone who didn" + [char]0x2019 + "t and didn" + [char]0x2019 + "t research the program enough could cause fatal errors to the system."

$issuesBlankInsertPoint = $issuesPara.Range.Duplicate
$issuesBlankInsertPoint.Start = $issuesPara.Range.End - 1
$issuesBlankInsertPoint.End = $issuesPara.Range.End - 1
$issuesBlankInsertPoint.InsertParagraphAfter()

# ---------------------------------------------------------------------------
# 6. Scalability paragraph: drop the spell-check markup around "Github" and
#    merge the three runs into one contiguous run of text.
# ---------------------------------------------------------------------------
$scalabilityHeading = Get-HeadingParagraph("Scalability")
$scalabilityBody = $scalabilityHeading.Next()

$searchRange = $scalabilityBody.Range.Duplicate
$found = $searchRange.Find.Execute("when multiple people are adding to the code and ")
if ($found) {
    $segStart = $searchRange.Start
    $segEnd = $scalabilityBody.Range.End - 1
    $placeholder = $d.Range($segStart, $segEnd)
    $placeholder.Text = [char]1
    $real = $d.Range($segStart, $segStart + 1)
    $real.Text = "when multiple people are adding to the code and Github commits can be completed more successfully. Also, due tot eh fact it is a command line based program, it is able to be added to slightly easier as new graphical interfaces don" + [char]0x2019 + "t have to be made and instead one can just create a new function and add a text based number code that locates it in the menu. "
}
